$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = "IMAGE_FAMILY"
$ws.Range("B10").Value = "IMAGE FAMILY"
$ws.Range("D10").Value = "URL_ATTRIBUTE,IMAGE_ATTRIBUTE,sku"
$ws.Range("G10").Value = "sku"
$ws.Range("H10").Value = "sku"
$ws.Range("I10").Value = "sku"

$ws.Columns.Item(6).ColumnWidth = 41.462962962963

$ws.Range("A16").Select()
